# Insert a new row at position 36, shifting rows 36-95 down to 37-96,
# then populate the new row 36 with the new data record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(36).Insert()

$ws.Cells.Item(36, 1).Value = 4
$ws.Cells.Item(36, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(36, 3).Value = "Los Lagos"
$ws.Cells.Item(36, 4).Value = 44536
$ws.Cells.Item(36, 5).Value = 10
$ws.Cells.Item(36, 6).Value = 100112022
$ws.Cells.Item(36, 7).Value = "Arveja Verde"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 35
$ws.Cells.Item(36, 11).Value = 20000
$ws.Cells.Item(36, 12).Value = 20000
$ws.Cells.Item(36, 13).Value = 20000
$ws.Cells.Item(36, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(36, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(36, 16).Value = 800
$ws.Cells.Item(36, 17).Value = 25
$ws.Cells.Item(36, 18).Value = "Hortaliza"
